# Weekly update: insert the newest week's two data rows (Repollo - Crespo record / Segunda)
# right after the header block of existing data (at row 1006), pushing the rest of the
# historical rows down by two. Because the very last two historical rows simply fall off
# the bottom of the originally-used range onto two brand-new rows, a plain row insert
# at 1006 reproduces the whole shift (including the two new trailing rows) for free.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 1006-1007; this shifts old rows 1006:1039 down to 1008:1041.
$ws.Range("A1006:A1007").EntireRow.Insert()

# New row 1006 - Repollo, Crespo record, Primera
$ws.Cells.Item(1006, 1).Value = 10
$ws.Cells.Item(1006, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1006, 3).Value = "La Araucanía"
$ws.Cells.Item(1006, 4).Value = 45075
$ws.Cells.Item(1006, 5).Value = 9
$ws.Cells.Item(1006, 6).Value = 100112006
$ws.Cells.Item(1006, 7).Value = "Repollo"
$ws.Cells.Item(1006, 8).Value = "Crespo record"
$ws.Cells.Item(1006, 9).Value = "Primera"
$ws.Cells.Item(1006, 10).Value = 1780
$ws.Cells.Item(1006, 11).Value = 1200
$ws.Cells.Item(1006, 12).Value = 1200
$ws.Cells.Item(1006, 13).Value = 1200
$ws.Cells.Item(1006, 14).Value = "$/unidad"
$ws.Cells.Item(1006, 15).Value = "Región del Maule"
$ws.Cells.Item(1006, 16).Value = 1200
$ws.Cells.Item(1006, 17).Value = 1
$ws.Cells.Item(1006, 18).Value = "Hortaliza"

# New row 1007 - Repollo, Crespo record, Segunda
$ws.Cells.Item(1007, 1).Value = 10
$ws.Cells.Item(1007, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1007, 3).Value = "La Araucanía"
$ws.Cells.Item(1007, 4).Value = 45075
$ws.Cells.Item(1007, 5).Value = 9
$ws.Cells.Item(1007, 6).Value = 100112006
$ws.Cells.Item(1007, 7).Value = "Repollo"
$ws.Cells.Item(1007, 8).Value = "Crespo record"
$ws.Cells.Item(1007, 9).Value = "Segunda"
$ws.Cells.Item(1007, 10).Value = 580
$ws.Cells.Item(1007, 11).Value = 1000
$ws.Cells.Item(1007, 12).Value = 1000
$ws.Cells.Item(1007, 13).Value = 1000
$ws.Cells.Item(1007, 14).Value = "$/unidad"
$ws.Cells.Item(1007, 15).Value = "Región del Maule"
$ws.Cells.Item(1007, 16).Value = 1000
$ws.Cells.Item(1007, 17).Value = 1
$ws.Cells.Item(1007, 18).Value = "Hortaliza"
